$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.013.97'
$ws.Range("E2").Value = '  +1.15%  '

# Row 3
$ws.Range("D3").Value = '1.759.36'
$ws.Range("E3").Value = '  +0.80%  '

# Row 4
$__style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").Style = $__style
$ws.Range("E4").Value = '  -0.44%  '

# Row 5
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.43'
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = '  -0.35%  '

# Row 6
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = '  -0.33%  '

# Row 7
$__style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5224'
$ws.Range("D7").Style = $__style
$ws.Range("E7").Value = '  +3.92%  '

# Row 8
$ws.Range("E8").Value = '  -3.05%  '

# Row 9
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2712'
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = '  +3.43%  '

# Row 10
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06208'
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = '  +1.17%  '

# Row 11
$ws.Range("D11").Value = '1.762.54'
$ws.Range("E11").Value = '  +0.78%  '

# Row 12
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07031'
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = '  +1.11%  '

# Row 13
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.76'
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = '  +3.62%  '

# Row 14
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6572'
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = '  +11.57%  '

# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.488'
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = '  +0.25%  '

# Row 16
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '78.24'
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = '  +1.98%  '

# Row 17
$__style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9988'
$ws.Range("D17").Style = $__style
$ws.Range("E17").Value = '  -0.48%  '

# Row 18
$__style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9995'
$ws.Range("D18").Style = $__style
$ws.Range("E18").Value = '  -0.36%  '

# Row 19
$ws.Range("D19").Value = '26.018.33'
$ws.Range("E19").Value = '  +0.94%  '

# Row 20
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.71'
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = '  +0.97%  '

# Row 21
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006710'
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = '  -1.05%  '

# Row 22
$ws.Range("D22").Value = '1.981.22'
$ws.Range("E22").Value = '  +0.49%  '

# Row 23
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.094'
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = '  +0.88%  '

# Row 24
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.426'
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = '  +4.14%  '

# Row 25
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.181'
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = '  +1.66%  '

# Row 26
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.20'
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = '  -0.83%  '

# Row 27
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.484'
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = '  -3.49%  '

# Row 28
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$__style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.835'
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = '  +1.21%  '

# Row 29
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.16'
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = '  +1.65%  '

# Row 30
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.65'
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = '  -0.50%  '

# Row 31
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08420'
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = '  +3.87%  '

# Row 32
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.703'
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = '  -1.53%  '

# Row 33
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.425'
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = '  -0.55%  '

# Row 34
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04423'
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = '  -1.57%  '

# Row 35
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.648'
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = '  +0.54%  '

# Row 36
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9992'
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = '  +2.46%  '

# Row 37
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6097'
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = '  +1.28%  '

# Row 38
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.736'
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = '  +3.16%  '

# Row 39
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01573'
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = '  +2.07%  '

# Row 40
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.958'
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = '  +2.84%  '

# Row 41
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = '  -0.16%  '

# Row 42
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.05'
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = '  -1.33%  '

# Row 43
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3895'
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = '  +3.22%  '

# Row 44
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7537'
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = '  +3.56%  '

# Row 45
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.936'
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = '  -3.75%  '

# Row 46
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05494'
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = '  +3.19%  '

# Row 47
$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1121'
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = '  +1.52%  '

# Row 48
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.126'
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = '  +4.27%  '

# Row 49
$ws.Range("E49").Value = '  +0.72%  '

# Row 50
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.69'
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = '  +0.86%  '

# Row 51
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = '  -0.02%  '
